$d = $word.ActiveDocument

# 1. Update the title text.
$d.Content.Find.Execute("Complex Test Document", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Test Document with Table", 2)

# 2. Update the intro paragraph text.
$d.Content.Find.Execute("This document has multiple tables.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is a test document.", 2)

# 3. Remove the first table entirely (Key/Value, Status/Active).
$d.Tables(1).Delete()

# 4. Remove the now-orphaned "Here is another table:" paragraph (including its mark).
$d.Paragraphs(3).Range.Delete()

# 5. Apply the LightGrid-Accent1 table style to the remaining table.
$t = $d.Tables(1)
$t.Style = "Light Grid Accent 1"

# 6. Update the header row text.
$t.Cell(1,1).Range.Text = "Name"
$t.Cell(1,2).Range.Text = "Age"
$t.Cell(1,3).Range.Text = "City"

# 7. Update row 2 (was Apple/1.99/100 -> Alice/30/NYC).
$t.Cell(2,1).Range.Text = "Alice"
$t.Cell(2,2).Range.Text = "30"
$t.Cell(2,3).Range.Text = "NYC"

# 8. Update row 3 (was Banana/0.99/50 -> Bob/25/LA).
$t.Cell(3,1).Range.Text = "Bob"
$t.Cell(3,2).Range.Text = "25"
$t.Cell(3,3).Range.Text = "LA"

# 9. Remove the last row (was Orange/2.49/75).
$t.Rows(4).Delete()
